# Apply the workbook fix described by the commit:
#   - row for IDDTE "09" (Dec-2032) on "Detalles" had the wrong tipo code -> should be "07"
#   - two new commission rows (Dec-2033 / Dec-2034) are appended on both
#     "Detalles" and "Apendices"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Detalles": fix row 11, append rows 12 and 13
# ---------------------------------------------------------------------
$detalles = $wb.Worksheets.Item("Detalles")

# Row 11 only changes its "IDDTE" code from 09 -> 07
$detalles.Cells.Item(11, 1).Value = "07"

# Helper fixed values shared by every detail row (columns B..H, J..N)
function Set-DetalleRow {
    param($sheet, $row, $codigo, $descripcion)

    $sheet.Cells.Item($row, 1).NumberFormat = "@"
    $sheet.Cells.Item($row, 1).Value = $codigo
    $sheet.Cells.Item($row, 2).Value = 1
    $sheet.Cells.Item($row, 3).Value = 2
    $sheet.Cells.Item($row, 4).Value = 1
    $sheet.Cells.Item($row, 5).Value = 20
    $sheet.Cells.Item($row, 6).Value = 99
    $sheet.Cells.Item($row, 7).Value = "null"
    $sheet.Cells.Item($row, 8).Value = "null"
    $sheet.Cells.Item($row, 9).Value = $descripcion
    $sheet.Cells.Item($row, 10).Value = 20
    $sheet.Cells.Item($row, 11).Value = 15
    $sheet.Cells.Item($row, 12).Value = 0
    $sheet.Cells.Item($row, 13).Value = 0
    $sheet.Cells.Item($row, 14).Value = 15
}

Set-DetalleRow $detalles 12 "03" "COMISION POR RECEPCION DE 0 PAGOS RECIBIDOS EN EL MES DE DICIEMBRE 2033"
Set-DetalleRow $detalles 13 "09" "COMISION POR RECEPCION DE 0 PAGOS RECIBIDOS EN EL MES DE DICIEMBRE 2034"

$detalles.Range("A13").Select()

# ---------------------------------------------------------------------
# Sheet "Apendices": append rows 12, 13 and 14
# ---------------------------------------------------------------------
$apendices = $wb.Worksheets.Item("Apendices")

function Set-ApendiceRow {
    param($sheet, $row, $codigo, $campo, $etiqueta, $valor)

    $sheet.Cells.Item($row, 1).NumberFormat = "@"
    $sheet.Cells.Item($row, 1).Value = $codigo
    $sheet.Cells.Item($row, 2).Value = $campo
    $sheet.Cells.Item($row, 3).Value = $etiqueta
    $sheet.Cells.Item($row, 4).Value = $valor
}

Set-ApendiceRow $apendices 12 "07" "Campo11" "Etiqueta11" "Valor11"
Set-ApendiceRow $apendices 13 "03" "Campo12" "Etiqueta12" "Valor12"
Set-ApendiceRow $apendices 14 "09" "Campo13" "Etiqueta13" "Valor13"

$apendices.Range("A14").Select()
